$wb = $excel.ActiveWorkbook

# Applies the numeric updates to the Leve profit-calculator tables,
# one per affected (sheet, row) pair, as produced by the scheduled
# market-price refresh run.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 834.93335
$ws.Range("I6").Value = 544.5714
$ws.Range("J6").Value = 4900
$ws.Range("K6").Value = 1633.7142
$ws.Range("L6").Value = 14700
$ws.Range("M6").Value = -1521.7142
$ws.Range("N6").Value = -14924

$ws.Range("H9").Value = 287.25
$ws.Range("I9").Value = 349.66666
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 349.66666
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = -180.66666
$ws.Range("N9").Value = -438

$ws.Range("H28").Value = 799.2308
$ws.Range("I28").Value = 642.3333
$ws.Range("K28").Value = 642.3333
$ws.Range("M28").Value = -157.3333

$ws.Range("H32").Value = 18135.5
$ws.Range("I32").Value = 13817.8
$ws.Range("J32").Value = 25331.666
$ws.Range("K32").Value = 13817.8
$ws.Range("L32").Value = 25331.666
$ws.Range("M32").Value = -13491.8
$ws.Range("N32").Value = -25983.666

$ws.Range("H38").Value = 1587
$ws.Range("I38").Value = 1699
$ws.Range("J38").Value = 1363
$ws.Range("K38").Value = 5097
$ws.Range("L38").Value = 4089
$ws.Range("M38").Value = -4725
$ws.Range("N38").Value = -4833

$ws.Range("H39").Value = 372.2857
$ws.Range("I39").Value = 391.9091
$ws.Range("J39").Value = 300.33334
$ws.Range("K39").Value = 1175.7273
$ws.Range("L39").Value = 901.0000200000001
$ws.Range("M39").Value = -879.7273
$ws.Range("N39").Value = -1493.00002

$ws.Range("H100").Value = 4157.8335
$ws.Range("I100").Value = 3210.5557
$ws.Range("K100").Value = 3210.5557
$ws.Range("M100").Value = -2669.5557

$ws.Range("H112").Value = 6114.923
$ws.Range("I112").Value = 949.5
$ws.Range("J112").Value = 6545.375
$ws.Range("K112").Value = 2848.5
$ws.Range("L112").Value = 19636.125
$ws.Range("M112").Value = -1740.5
$ws.Range("N112").Value = -21852.125

$ws.Range("H129").Value = 2124.5
$ws.Range("J129").Value = 5000
$ws.Range("L129").Value = 15000
$ws.Range("N129").Value = -25000

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 1896.6666
$ws.Range("I26").Value = 1896.6666
$ws.Range("K26").Value = 1896.6666
$ws.Range("M26").Value = -1566.6666

$ws.Range("H32").Value = 21440.371
$ws.Range("I32").Value = 22795.164
$ws.Range("K32").Value = 22795.164
$ws.Range("M32").Value = -22508.164

$ws.Range("H61").Value = 5981.0713
$ws.Range("I61").Value = 3810.3438
$ws.Range("K61").Value = 3810.3438
$ws.Range("M61").Value = -3598.3438

$ws.Range("H136").Value = 5981.0713
$ws.Range("I136").Value = 3810.3438
$ws.Range("K136").Value = 11431.0314
$ws.Range("M136").Value = -8881.0314

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 82759.8
$ws.Range("I20").Value = 78449.75
$ws.Range("J20").Value = 100000
$ws.Range("K20").Value = 78449.75
$ws.Range("L20").Value = 100000
$ws.Range("M20").Value = -78202.75
$ws.Range("N20").Value = -100494

$ws.Range("H76").Value = 29413.375
$ws.Range("J76").Value = 29413.375
$ws.Range("L76").Value = 29413.375
$ws.Range("N76").Value = -30043.375

$ws.Range("H79").Value = 29413.375
$ws.Range("J79").Value = 29413.375
$ws.Range("L79").Value = 29413.375
$ws.Range("N79").Value = -31597.375

$ws.Range("H80").Value = 572.2778
$ws.Range("I80").Value = 896.125
$ws.Range("K80").Value = 896.125
$ws.Range("M80").Value = 101.875

$ws.Range("H82").Value = 17023.928
$ws.Range("I82").Value = 9833.5
$ws.Range("K82").Value = 9833.5
$ws.Range("M82").Value = -9450.5

$ws.Range("H83").Value = 572.2778
$ws.Range("I83").Value = 896.125
$ws.Range("K83").Value = 4480.625
$ws.Range("M83").Value = 511.375

$ws.Range("H85").Value = 17023.928
$ws.Range("I85").Value = 9833.5
$ws.Range("K85").Value = 9833.5
$ws.Range("M85").Value = -8507.5

$ws.Range("H86").Value = 2436.3333
$ws.Range("I86").Value = 2549.6428
$ws.Range("K86").Value = 2549.6428
$ws.Range("M86").Value = -1426.6428

$ws.Range("H89").Value = 2436.3333
$ws.Range("I89").Value = 2549.6428
$ws.Range("K89").Value = 12748.214
$ws.Range("M89").Value = -7132.214

$ws.Range("H94").Value = 5883130
$ws.Range("I94").Value = 853.61536
$ws.Range("J94").Value = 25000528
$ws.Range("K94").Value = 853.61536
$ws.Range("L94").Value = 25000528
$ws.Range("M94").Value = -402.61536
$ws.Range("N94").Value = -25001430

$ws.Range("H105").Value = 5095.7
$ws.Range("I105").Value = 7941
$ws.Range("K105").Value = 7941
$ws.Range("M105").Value = -6194

$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2383486.8
$ws.Range("I31").Value = 5002498.5
$ws.Range("J31").Value = 2567.1365
$ws.Range("K31").Value = 5002498.5
$ws.Range("L31").Value = 2567.1365
$ws.Range("M31").Value = -5002203.5
$ws.Range("N31").Value = -3157.1365

$ws.Range("H32").Value = 40000
$ws.Range("I32").Value = 40000
$ws.Range("K32").Value = 40000
$ws.Range("M32").Value = -39684

$ws.Range("H34").Value = 2383486.8
$ws.Range("I34").Value = 5002498.5
$ws.Range("J34").Value = 2567.1365
$ws.Range("K34").Value = 5002498.5
$ws.Range("L34").Value = 2567.1365
$ws.Range("M34").Value = -5002296.5
$ws.Range("N34").Value = -2971.1365

$ws.Range("H99").Value = 4475.125
$ws.Range("I99").Value = 4473.75
$ws.Range("J99").Value = 4475.5835
$ws.Range("K99").Value = 4473.75
$ws.Range("L99").Value = 4475.5835
$ws.Range("M99").Value = -2975.75
$ws.Range("N99").Value = -7471.5835

$ws.Range("H126").Value = 4475.125
$ws.Range("I126").Value = 4473.75
$ws.Range("J126").Value = 4475.5835
$ws.Range("K126").Value = 13421.25
$ws.Range("L126").Value = 13426.7505
$ws.Range("M126").Value = -10951.25
$ws.Range("N126").Value = -18366.7505

$ws.Range("H134").Value = 2388.7896
$ws.Range("J134").Value = 4256.75
$ws.Range("L134").Value = 12770.25
$ws.Range("N134").Value = -17840.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1498.3334
$ws.Range("I92").Value = 1733
$ws.Range("J92").Value = 794.3333
$ws.Range("K92").Value = 5199
$ws.Range("L92").Value = 2382.9999
$ws.Range("M92").Value = -3951
$ws.Range("N92").Value = -4878.9999

$ws.Range("H98").Value = 350
$ws.Range("I98").Value = 100
$ws.Range("K98").Value = 300
$ws.Range("M98").Value = 1198

$ws.Range("H105").Value = 7353.25
$ws.Range("J105").Value = 7685.7144
$ws.Range("L105").Value = 23057.1432
$ws.Range("N105").Value = -28299.1432

$ws.Range("H140").Value = 72593
$ws.Range("I140").Value = 72593
$ws.Range("K140").Value = 217779
$ws.Range("M140").Value = -212599

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 55557780
$ws.Range("I122").Value = 1081
$ws.Range("J122").Value = 83336130
$ws.Range("K122").Value = 3243
$ws.Range("L122").Value = 250008390
$ws.Range("M122").Value = -793
$ws.Range("N122").Value = -250013290

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1075.5625
$ws.Range("I16").Value = 979.0714
$ws.Range("K16").Value = 979.0714
$ws.Range("M16").Value = -809.0714

$ws.Range("H55").Value = 1950.875
$ws.Range("I55").Value = 689
$ws.Range("K55").Value = 689
$ws.Range("M55").Value = -516

$ws.Range("H82").Value = 4305.2856
$ws.Range("I82").Value = 4111.385
$ws.Range("J82").Value = 4620.375
$ws.Range("K82").Value = 4111.385
$ws.Range("L82").Value = 4620.375
$ws.Range("M82").Value = -3750.385
$ws.Range("N82").Value = -5342.375

$ws.Range("H85").Value = 4305.2856
$ws.Range("I85").Value = 4111.385
$ws.Range("J85").Value = 4620.375
$ws.Range("K85").Value = 4111.385
$ws.Range("L85").Value = 4620.375
$ws.Range("M85").Value = -2863.385
$ws.Range("N85").Value = -7116.375

$ws.Range("H93").Value = 1061.2222
$ws.Range("I93").Value = 1135.7142
$ws.Range("J93").Value = 800.5
$ws.Range("K93").Value = 1135.7142
$ws.Range("L93").Value = 800.5
$ws.Range("M93").Value = 112.2858000000001
$ws.Range("N93").Value = -3296.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4803.9644
$ws.Range("I81").Value = 4713.7827
$ws.Range("K81").Value = 9427.565399999999
$ws.Range("M81").Value = -8366.565399999999

$ws.Range("H84").Value = 4803.9644
$ws.Range("I84").Value = 4713.7827
$ws.Range("K84").Value = 47137.827
$ws.Range("M84").Value = -41833.827

$ws.Range("H100").Value = 1351.091
$ws.Range("I100").Value = 896.3333
$ws.Range("J100").Value = 1896.8
$ws.Range("K100").Value = 1792.6666
$ws.Range("L100").Value = 3793.6
$ws.Range("M100").Value = -1251.6666
$ws.Range("N100").Value = -4875.6

$ws.Range("H113").Value = 512.82855
$ws.Range("I113").Value = 582.1539
$ws.Range("K113").Value = 1746.4617
$ws.Range("M113").Value = 423.5382999999999

$ws.Range("H126").Value = 2871.6667
$ws.Range("I126").Value = 2308.923
$ws.Range("K126").Value = 6926.768999999999
$ws.Range("M126").Value = -4456.768999999999

$ws.Range("H136").Value = 7934.636
$ws.Range("I136").Value = 11972
$ws.Range("K136").Value = 35916
$ws.Range("M136").Value = -33366

